$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 02:22"

# Estados Unidos (row 4) - new daily totals
$ws.Range("B4").Value = 764303
$ws.Range("C4").Value = 25511
$ws.Range("D4").Value = 71003
$ws.Range("E4").Value = 652752
$ws.Range("F4").Value = 13566
$ws.Range("G4").Value = 1534
$ws.Range("H4").Value = 40548

# Alemania (row 8) - new daily totals
$ws.Range("B8").Value = 145742
$ws.Range("C8").Value = 2018
$ws.Range("D8").Value = 88000
$ws.Range("E8").Value = 53100
$ws.Range("F8").Value = 2889
$ws.Range("G8").Value = 104
$ws.Range("H8").Value = 4642

# Argentina overtakes Marruecos -> rows 55/56 swap order and get updated data
$ws.Range("A55").Value = "Argentina"
$ws.Range("B55").Value = 2941
$ws.Range("C55").Value = 102
$ws.Range("D55").Value = 709
$ws.Range("E55").Value = 2098
$ws.Range("F55").Value = 123
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 134

$ws.Range("A56").Value = "Marruecos"
$ws.Range("B56").Value = 2855
$ws.Range("C56").Value = 170
$ws.Range("D56").Value = 327
$ws.Range("E56").Value = 2387
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 141

# Bahamas overtakes Polinesia Francesa -> rows 156/157 swap order and get updated data
$ws.Range("A156").Value = "Bahamas"
$ws.Range("B156").Value = 60
$ws.Range("C156").Value = 2
$ws.Range("D156").Value = 11
$ws.Range("E156").Value = 40
$ws.Range("F156").Value = 1
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 9

$ws.Range("A157").Value = "Polinesia Francesa"
$ws.Range("B157").Value = 55
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 2
$ws.Range("E157").Value = 53
$ws.Range("F157").Value = 1
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 0

# Mozambique overtakes Siria -> rows 166/167 swap order and get updated data
$ws.Range("A166").Value = "Mozambique"
$ws.Range("B166").Value = 39
$ws.Range("C166").Value = 4
$ws.Range("D166").Value = 8
$ws.Range("E166").Value = 31
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 0

$ws.Range("A167").Value = "Siria"
$ws.Range("B167").Value = 39
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 5
$ws.Range("E167").Value = 31
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 3
